$wb = $excel.ActiveWorkbook

# --- PredicateObjectMaps sheet: data + formatting + hyperlink updates ---
$ws = $wb.Worksheets.Item("PredicateObjectMaps")

# Row 4 (schema:description / EMPLOYED): Object becomes a URI, DataType becomes anyURI
$ws.Range("B4").Value = "http://ex.com/{bio}"
$ws.Range("C4").Value = "anyURI"

# Row 11 (schema:articleBody / NEW): DataType becomes iri (trailing space kept)
$ws.Range("C11").Value = "iri "

# Normalize formatting on rows 10-11 to match the rest of the table (style 2)
$ws.Range("A2").Copy()
$ws.Range("A10:B11").PasteSpecial(-4122)

# Add hyperlink on B4 pointing at the new URI, then restore its formatting
# (Hyperlinks.Add applies Excel's default hyperlink style; paste the original
# format back over it so the cell keeps looking like its neighbours)
$h = $ws.Hyperlinks.Add($ws.Range("B4"), "http://ex.com/{bio}")
$h.TextToDisplay = "http://ex.com/{bio}"
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Selection / active cell on this sheet
$ws.Range("C12").Select()

# --- Active sheet / tab selection bookkeeping ---
$wsFunctions = $wb.Worksheets.Item("Functions")
$wsFunctions.Select()
$ws.Select()
